$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (UsedRange may start at row 0)
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Column C ("Förändrad") holds a date serial number (45188 = 2023-09-19) for
# every data row (rows 2..lastRow). Bump every one of them by one day
# (45188 -> 45189), leaving all other columns untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
